$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 15 closes out its script-group: re-style it like the other "last
#    row of a group" rows (3, 6, 10, 12) which carry a bottom border
#    (style 5 for columns A/B, style 6 for columns C/D/E) instead of the
#    plain wrap style (3 / 4). Copying the formats from row 3 reuses the
#    existing style entries instead of inventing new ones.
# ---------------------------------------------------------------------------
$ws.Range("A3:E3").Copy()
$ws.Range("A15:E15").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Fill in the new data for rows 16-18 (a new script entry,
#    SCRIPT/G01P03A/us0105.ssb). Values are entered in the same column-major
#    order the original author used, so new shared-string entries line up
#    the same way: English column first, then the filename, then the
#    Russian translation column, then the encoded/ciphered column.
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = " For some reason, no sunlight\nfalls on [CS:P]Luminous Spring[CR]."
$ws.Range("C17").Value = " ...[K]Huh? Why\'s that?"
$ws.Range("C18").Value = " Well, why would I know that?"

$ws.Range("A16").Value = "SCRIPT/G01P03A/us0105.ssb"

$ws.Range("D16").Value = " По какой-то причине, на [CS:P]Сияющий\nИсточник[CR] не светит солнце."
$ws.Range("D17").Value = " ...[K]Что? Почему?"
$ws.Range("D18").Value = " Ну, а я то откуда знаю?"

$ws.Range("E16").Value = " Ðï ëàëïê-óï ðñéœéîå, îà [CS:P]Òéÿýþéê\nÉòóïœîéë[CR] îå òâåóéó òïìîøå."
$ws.Range("E17").Value = " ...[K]Œóï? Ðïœåíô?"
$ws.Range("E18").Value = " Îô, à ÿ óï ïóëôäà èîàý?"

$ws.Range("B16").Value = 57
$ws.Range("B17").Value = 60
$ws.Range("B18").Value = 63

# ---------------------------------------------------------------------------
# 3) Style the new rows like the other "interior group" rows (2, 4, 7 ...):
#    plain wrap style, no border - style 3 for columns A/B, style 4 for C/D/E.
# ---------------------------------------------------------------------------
$ws.Range("A2:E2").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B5:E5").Copy()
$ws.Range("B17:E17").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B18:E18").PasteSpecial(-4122)   # xlPasteFormats

# Row 16 keeps the taller, multi-line auto height like the other 3-line rows.
$ws.Rows.Item(16).RowHeight = 43.2

# ---------------------------------------------------------------------------
# 4) Move the active selection the way the author left it (D15, scrolled down
#    a bit further than before).
# ---------------------------------------------------------------------------
$ws.Range("A13:E18").Select()
$ws.Range("D15").Activate()
